# Update the "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" sheets, reflecting the newly generated output.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# row number -> new value for column F
$updates = @{
    2  = 162
    8  = 12104
    15 = 13528
    16 = 13580
    25 = 189
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
